# Sprint 15 (pagination + bug fixing) edits to the "Reports" sheet:
#   - update the organisation name, ledger-generation date, and report type
#   - drop the stray CREDITNOTE detail row + its TOTAL row (pagination clean-up)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "fkjds" -> "rk industries"
$ws.Range("B2").Value = "rk industries"

# B5: "2020-08-11" -> "2020-09-11"
# Setting .Value/.Formula directly with an ISO-looking date string makes
# Excel auto-convert it to a real date serial, which isn't what we want here
# (the source data is a plain text date label). Write it as a text formula
# first, then flatten the cell to a static value via copy / paste-values so
# the stored cell is the literal text "2020-09-11" with no lingering formula.
$ws.Range("B5").Formula = '="2020-09-11"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

# B6: "Yearly" -> "Monthly"
$ws.Range("B6").Value = "Monthly"

# Remove the "new with tax" CREDITNOTE row (10) and the TOTAL row below it
# (11) entirely -- not just clear their contents -- so the sheet's used
# range shrinks back down to A1:E9. Delete bottom-up so the row indices
# don't shift out from under the second call.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
